# Insert a new price-report row for "Región de Coquimbo" (date 44740) right
# before the existing row 45 (Feria Lagunitas de Puerto Montt - Haba), which
# pushes that row and every row after it down by one (old row 85 becomes 86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45..85 down to 46..86, freeing up row 45 for the new record.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(45, 1).Value  = 4
$ws.Cells.Item(45, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value  = "Los Lagos"
$ws.Cells.Item(45, 4).Value  = 44740
$ws.Cells.Item(45, 5).Value  = 10
$ws.Cells.Item(45, 6).Value  = 100112026
$ws.Cells.Item(45, 7).Value  = "Haba"
$ws.Cells.Item(45, 8).Value  = "Sin especificar"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 80
$ws.Cells.Item(45, 11).Value = 23000
$ws.Cells.Item(45, 12).Value = 23000
$ws.Cells.Item(45, 13).Value = 23000
$ws.Cells.Item(45, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(45, 16).Value = 920
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
